$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 9445.666999999999
$ws.Range("I55").Value = 2690
$ws.Range("J55").Value = 12823.5
$ws.Range("K55").Value = 2690
$ws.Range("L55").Value = 12823.5
$ws.Range("M55").Value = -2476
$ws.Range("N55").Value = -13251.5
$ws.Range("H129").Value = 845.5806
$ws.Range("J129").Value = 920.3461
$ws.Range("L129").Value = 2761.0383
$ws.Range("N129").Value = -12761.0383
$ws.Range("H132").Value = 34506.535
$ws.Range("I132").Value = 44408
$ws.Range("J132").Value = 1973.1428
$ws.Range("K132").Value = 133224
$ws.Range("L132").Value = 5919.428400000001
$ws.Range("M132").Value = -130694
$ws.Range("N132").Value = -10979.4284
$ws.Range("H138").Value = 3179.3396
$ws.Range("J138").Value = 3690.2974
$ws.Range("L138").Value = 11070.8922
$ws.Range("N138").Value = -21350.8922
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3478.7778
$ws.Range("I132").Value = 1740.1
$ws.Range("J132").Value = 4501.5293
$ws.Range("K132").Value = 5220.299999999999
$ws.Range("L132").Value = 13504.5879
$ws.Range("M132").Value = -2690.299999999999
$ws.Range("N132").Value = -18564.5879
$ws.Range("H133").Value = 59230.5
$ws.Range("J133").Value = 59230.5
$ws.Range("L133").Value = 59230.5
$ws.Range("N133").Value = -64290.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H105").Value = 2395.5173
$ws.Range("I105").Value = 2006.36
$ws.Range("J105").Value = 4827.75
$ws.Range("K105").Value = 2006.36
$ws.Range("L105").Value = 4827.75
$ws.Range("M105").Value = -259.3599999999999
$ws.Range("N105").Value = -8321.75
$ws.Range("H134").Value = 2441.1936
$ws.Range("I134").Value = 1211.8572
$ws.Range("J134").Value = 5022.8
$ws.Range("K134").Value = 3635.5716
$ws.Range("L134").Value = 15068.4
$ws.Range("M134").Value = -1100.5716
$ws.Range("N134").Value = -20138.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 190
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 235
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 235
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = -935
$ws.Range("H99").Value = 6559.579
$ws.Range("I99").Value = 1318.1666
$ws.Range("K99").Value = 1318.1666
$ws.Range("M99").Value = 179.8334
$ws.Range("H126").Value = 6559.579
$ws.Range("I126").Value = 1318.1666
$ws.Range("K126").Value = 3954.4998
$ws.Range("M126").Value = -1484.4998
$ws.Range("H135").Value = 44078
$ws.Range("J135").Value = 44078
$ws.Range("L135").Value = 44078
$ws.Range("N135").Value = -54218
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1432.421
$ws.Range("I5").Value = 729.7143
$ws.Range("J5").Value = 3400
$ws.Range("K5").Value = 2189.1429
$ws.Range("L5").Value = 10200
$ws.Range("M5").Value = -2077.1429
$ws.Range("N5").Value = -10424
$ws.Range("H8").Value = 119.888885
$ws.Range("I8").Value = 119.888885
$ws.Range("K8").Value = 359.666655
$ws.Range("M8").Value = -220.666655
$ws.Range("H9").Value = 260000
$ws.Range("J9").Value = 260000
$ws.Range("L9").Value = 780000
$ws.Range("N9").Value = -780448
$ws.Range("H98").Value = 650.25
$ws.Range("I98").Value = 1354
$ws.Range("K98").Value = 4062
$ws.Range("M98").Value = -2564
$ws.Range("H135").Value = 1432.421
$ws.Range("I135").Value = 729.7143
$ws.Range("J135").Value = 3400
$ws.Range("K135").Value = 6567.428699999999
$ws.Range("L135").Value = 30600
$ws.Range("M135").Value = -4032.428699999999
$ws.Range("N135").Value = -35670
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2795.6428
$ws.Range("I80").Value = 2626.6667
$ws.Range("K80").Value = 2626.6667
$ws.Range("M80").Value = -1628.6667
$ws.Range("H83").Value = 2795.6428
$ws.Range("I83").Value = 2626.6667
$ws.Range("K83").Value = 13133.3335
$ws.Range("M83").Value = -8141.333500000001
$ws.Range("H93").Value = 19800
$ws.Range("J93").Value = 19800
$ws.Range("L93").Value = 19800
$ws.Range("N93").Value = -23544
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H133").Value = 30695
$ws.Range("J133").Value = 30695
$ws.Range("L133").Value = 30695
$ws.Range("N133").Value = -35755
$ws.Range("H136").Value = 8774212
$ws.Range("I136").Value = 1690.9678
$ws.Range("K136").Value = 5072.903399999999
$ws.Range("M136").Value = -2522.903399999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 8714
$ws.Range("I41").Value = 6421
$ws.Range("J41").Value = 13300
$ws.Range("K41").Value = 6421
$ws.Range("L41").Value = 13300
$ws.Range("M41").Value = -6031
$ws.Range("N41").Value = -14080
$ws.Range("H45").Value = 7750
$ws.Range("J45").Value = 7750
$ws.Range("L45").Value = 7750
$ws.Range("N45").Value = -8732
$ws.Range("H81").Value = 882.45
$ws.Range("I81").Value = 726
$ws.Range("J81").Value = 1038.9
$ws.Range("K81").Value = 1452
$ws.Range("L81").Value = 2077.8
$ws.Range("M81").Value = -391
$ws.Range("N81").Value = -4199.8
$ws.Range("H84").Value = 882.45
$ws.Range("I84").Value = 726
$ws.Range("J84").Value = 1038.9
$ws.Range("K84").Value = 7260
$ws.Range("L84").Value = 10389
$ws.Range("M84").Value = -1956
$ws.Range("N84").Value = -20997
$ws.Range("H86").Value = 25291.666
$ws.Range("J86").Value = 25291.666
$ws.Range("L86").Value = 25291.666
$ws.Range("N86").Value = -27537.666
$ws.Range("H89").Value = 25291.666
$ws.Range("J89").Value = 25291.666
$ws.Range("L89").Value = 126458.33
$ws.Range("N89").Value = -137690.33
$ws.Range("H96").Value = 1751.25
$ws.Range("I96").Value = 1841.2
$ws.Range("K96").Value = 1841.2
$ws.Range("M96").Value = -468.2
$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982
$ws.Range("H100").Value = 542.8570999999999
$ws.Range("I100").Value = 466.66666
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 933.33332
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -392.33332
$ws.Range("N100").Value = -3082
$ws.Range("H122").Value = 34533.484
$ws.Range("I122").Value = 73666.71000000001
$ws.Range("J122").Value = 2306.1177
$ws.Range("K122").Value = 221000.13
$ws.Range("L122").Value = 6918.353099999999
$ws.Range("M122").Value = -218550.13
$ws.Range("N122").Value = -11818.3531
$ws.Range("H136").Value = 2550.9143
$ws.Range("I136").Value = 951.06665
$ws.Range("J136").Value = 3750.8
$ws.Range("K136").Value = 2853.19995
$ws.Range("L136").Value = 11252.4
$ws.Range("M136").Value = -303.1999500000002
$ws.Range("N136").Value = -16352.4
